# Little update of presentation
$p = $ppt.ActivePresentation

# Remove the second slide (the picture-only slide).
$p.Slides.Item(2).Delete()

# Add a title textbox to the remaining slide with "Cities of Beacons".
$s = $p.Slides.Item(1)
$box = $s.Shapes.AddTextbox(1, 0, 87.21125984251968, 960, 50.892204724409446)
$box.Name = "TextovéPole 1"
$box.Fill.Visible = $false

$tf = $box.TextFrame
$tf.WordWrap = $true
$tf.AutoSize = 1

$tr = $tf.TextRange
$tr.Text = "Cities"
$tr.Font.NameComplexScript = "Times New Roman"

$r2 = $tr.InsertAfter(" ")
$r2.Font.NameComplexScript = "Times New Roman"

$r3 = $tr.InsertAfter("of")
$r3.Font.NameComplexScript = "Times New Roman"

$r4 = $tr.InsertAfter(" ")
$r4.Font.NameComplexScript = "Times New Roman"

$r5 = $tr.InsertAfter("Beacons")
$r5.Font.NameComplexScript = "Times New Roman"

$tr.ParagraphFormat.Alignment = 2
$tr.Font.Name = "Times New Roman"
$tr.Font.Size = 36
